$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revised AgTests (H) / AgPosit (I) figures for existing rows (antigen test data correction)
$ws.Range("H268").Value = 13814
$ws.Range("I268").Value = 666
$ws.Range("H271").Value = 42725
$ws.Range("I271").Value = 1655
$ws.Range("H272").Value = 30441
$ws.Range("I272").Value = 1638
$ws.Range("H273").Value = 26934
$ws.Range("I273").Value = 1343
$ws.Range("H274").Value = 28546
$ws.Range("I274").Value = 1325
$ws.Range("H275").Value = 29139
$ws.Range("I275").Value = 1251
$ws.Range("H276").Value = 11741
$ws.Range("I276").Value = 418
$ws.Range("H278").Value = 30008
$ws.Range("I278").Value = 2128
$ws.Range("H279").Value = 42946
$ws.Range("I279").Value = 3100
$ws.Range("H280").Value = 35112
$ws.Range("I280").Value = 2351
$ws.Range("H281").Value = 45427
$ws.Range("I281").Value = 3286
$ws.Range("H282").Value = 47241
$ws.Range("I282").Value = 2845
$ws.Range("H283").Value = 16913
$ws.Range("I283").Value = 1001
$ws.Range("H285").Value = 41494
$ws.Range("I285").Value = 3462
$ws.Range("H286").Value = 54662
$ws.Range("I286").Value = 4267
$ws.Range("H287").Value = 57746
$ws.Range("I287").Value = 3920
$ws.Range("H288").Value = 57099
$ws.Range("I288").Value = 4000
$ws.Range("H289").Value = 63138
$ws.Range("I289").Value = 3623
$ws.Range("H290").Value = 17626
$ws.Range("I290").Value = 1474
$ws.Range("H291").Value = 14851
$ws.Range("H292").Value = 81769
$ws.Range("I292").Value = 7255
$ws.Range("H293").Value = 81657
$ws.Range("I293").Value = 5745
$ws.Range("H294").Value = 90992
$ws.Range("I294").Value = 5040
$ws.Range("H295").Value = 17958
$ws.Range("I295").Value = 1056
$ws.Range("H299").Value = 64544
$ws.Range("I299").Value = 6802
$ws.Range("H300").Value = 71397
$ws.Range("I300").Value = 7028
$ws.Range("H301").Value = 69695
$ws.Range("I301").Value = 5555
$ws.Range("H302").Value = 74110
$ws.Range("I302").Value = 5364
$ws.Range("H303").Value = 10165
$ws.Range("I303").Value = 665
$ws.Range("H304").Value = 6511
$ws.Range("H306").Value = 71105
$ws.Range("I306").Value = 7209
$ws.Range("H307").Value = 74660
$ws.Range("I307").Value = 6415
$ws.Range("H308").Value = 15688
$ws.Range("I308").Value = 1327
$ws.Range("H309").Value = 58840
$ws.Range("I309").Value = 4056
$ws.Range("H310").Value = 91546
$ws.Range("I310").Value = 5209
$ws.Range("H311").Value = 36065
$ws.Range("I311").Value = 1335
$ws.Range("H312").Value = 40279
$ws.Range("H313").Value = 73967
$ws.Range("I313").Value = 3596
$ws.Range("H314").Value = 64641
$ws.Range("I314").Value = 3318
$ws.Range("H315").Value = 65879
$ws.Range("I315").Value = 2723
$ws.Range("H316").Value = 48644
$ws.Range("I316").Value = 2223
$ws.Range("H317").Value = 61584
$ws.Range("I317").Value = 2123
$ws.Range("H318").Value = 23054
$ws.Range("I318").Value = 872
$ws.Range("H319").Value = 56054
$ws.Range("I319").Value = 1778
$ws.Range("H320").Value = 87122
$ws.Range("I320").Value = 3942
$ws.Range("H321").Value = 90259
$ws.Range("I321").Value = 2789
$ws.Range("H322").Value = 105138
$ws.Range("I322").Value = 2272
$ws.Range("H323").Value = 147599
$ws.Range("I323").Value = 2289
$ws.Range("H324").Value = 230408
$ws.Range("I324").Value = 2668
$ws.Range("H325").Value = 684258
$ws.Range("I325").Value = 5624
$ws.Range("H326").Value = 408272
$ws.Range("I326").Value = 3601
$ws.Range("H327").Value = 257844
$ws.Range("I327").Value = 3689
$ws.Range("H328").Value = 187092
$ws.Range("I328").Value = 2725
$ws.Range("H329").Value = 82182
$ws.Range("I329").Value = 1867
$ws.Range("H330").Value = 68480
$ws.Range("I330").Value = 1961

# Append new daily rows for 2021-01-29, 2021-01-30, 2021-01-31
$ws.Range("A331").Value = 44225
$ws.Range("B331").Value = 248190
$ws.Range("C331").Value = 212271
$ws.Range("D331").Value = 31354
$ws.Range("E331").Value = 11207
$ws.Range("F331").Value = 2182
$ws.Range("G331").Value = 4565
$ws.Range("H331").Value = 142646
$ws.Range("I331").Value = 2462
$ws.Range("A332").Value = 44226
$ws.Range("B332").Value = 249913
$ws.Range("C332").Value = 216052
$ws.Range("D332").Value = 29219
$ws.Range("E332").Value = 8201
$ws.Range("F332").Value = 1723
$ws.Range("G332").Value = 4642
$ws.Range("H332").Value = 349642
$ws.Range("I332").Value = 4022
$ws.Range("A333").Value = 44227
$ws.Range("B333").Value = 250357
$ws.Range("C333").Value = 218923
$ws.Range("D333").Value = 26723
$ws.Range("E333").Value = 2290
$ws.Range("F333").Value = 444
$ws.Range("G333").Value = 4711
$ws.Range("H333").Value = 204733
$ws.Range("I333").Value = 2127
